$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-90 down to 45-91
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new data
$ws.Cells.Item(44, 1).Value = 9
$ws.Cells.Item(44, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(44, 3).Value = "Metropolitana"
$ws.Cells.Item(44, 4).Value = 44827
$ws.Cells.Item(44, 5).Value = 13
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100102
$ws.Cells.Item(44, 8).Value = "Cítricos"
$ws.Cells.Item(44, 9).Value = 100102006
$ws.Cells.Item(44, 10).Value = "Pomelo"
$ws.Cells.Item(44, 11).Value = "Start Ruby"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 150
$ws.Cells.Item(44, 14).Value = 12000
$ws.Cells.Item(44, 15).Value = 12000
$ws.Cells.Item(44, 16).Value = 12000
$ws.Cells.Item(44, 17).Value = "`$/caja 14 kilos"
$ws.Cells.Item(44, 18).Value = "Región Metropolitana"
$ws.Cells.Item(44, 19).Value = 857
$ws.Cells.Item(44, 20).Value = 14
